$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Table_4_15")

# Update report title (October 2016 -> November 2016)
$ws.Range("A1").Value = "Table 4.15. Receipts and Quality of Coal by Rank Delivered for Electricity Generation: Electric Utilties by State, November 2016"

# Updated data cells (EPM November-2016 preliminary run)
$ws.Range("B5").Value = 0
$ws.Range("C5").Value = "--"
$ws.Range("D5").Value = "--"
$ws.Range("B9").Value = 0
$ws.Range("C9").Value = "--"
$ws.Range("D9").Value = "--"
$ws.Range("B16").Value = 3054
$ws.Range("C16").Value = 2.89
$ws.Range("E16").Value = 3697
$ws.Range("G16").Value = 5.0999999999999996
$ws.Range("B17").Value = 144
$ws.Range("C17").Value = 2.8
$ws.Range("D17").Value = 12.7
$ws.Range("E17").Value = 412
$ws.Range("F17").Value = 0.2
$ws.Range("G17").Value = 4.5999999999999996
$ws.Range("B18").Value = 2099
$ws.Range("C18").Value = 2.74
$ws.Range("D18").Value = 8.5
$ws.Range("E18").Value = 117
$ws.Range("G18").Value = 4.7
$ws.Range("B19").Value = 84
$ws.Range("C19").Value = 2.67
$ws.Range("E19").Value = 1487
$ws.Range("G19").Value = 5
$ws.Range("B20").Value = 589
$ws.Range("C20").Value = 3.45
$ws.Range("D20").Value = 8.5
$ws.Range("B21").Value = 137
$ws.Range("C21").Value = 2.87
$ws.Range("E21").Value = 1681
$ws.Range("G21").Value = 5.3
$ws.Range("C22").Value = 3.08
$ws.Range("E22").Value = 7912
$ws.Range("F22").Value = 0.28999999999999998
$ws.Range("H22").Value = 1755
$ws.Range("I22").Value = 0.79
$ws.Range("J22").Value = 10.4
$ws.Range("E23").Value = 1291
$ws.Range("F23").Value = 0.24
$ws.Range("C24").Value = 3.08
$ws.Range("E24").Value = 1303
$ws.Range("F24").Value = 0.31
$ws.Range("G24").Value = 4.8
$ws.Range("E25").Value = 1395
$ws.Range("F25").Value = 0.41
$ws.Range("G25").Value = 6.3
$ws.Range("E26").Value = 2713
$ws.Range("F26").Value = 0.23
$ws.Range("E27").Value = 1015
$ws.Range("F27").Value = 0.28999999999999998
$ws.Range("E28").Value = 39
$ws.Range("G28").Value = 5.3
$ws.Range("H28").Value = 1755
$ws.Range("I28").Value = 0.79
$ws.Range("J28").Value = 10.4
$ws.Range("E29").Value = 156
$ws.Range("F29").Value = 0.35
$ws.Range("G29").Value = 5.2
$ws.Range("B30").Value = 6019
$ws.Range("C30").Value = 2.37
$ws.Range("D30").Value = 10
$ws.Range("E30").Value = 843
$ws.Range("G30").Value = 4.8
$ws.Range("B33").Value = 1169
$ws.Range("C33").Value = 2.38
$ws.Range("D33").Value = 8.6
$ws.Range("B34").Value = 690
$ws.Range("C34").Value = 2.56
$ws.Range("D34").Value = 7.7
$ws.Range("E34").Value = 843
$ws.Range("G34").Value = 4.8
$ws.Range("B36").Value = 1205
$ws.Range("C36").Value = 1.94
$ws.Range("D36").Value = 9.5
$ws.Range("B37").Value = 632
$ws.Range("C37").Value = 1.76
$ws.Range("D37").Value = 9.1999999999999993
$ws.Range("B38").Value = 447
$ws.Range("C38").Value = 1.1100000000000001
$ws.Range("D38").Value = 15.5
$ws.Range("B39").Value = 1877
$ws.Range("C39").Value = 3.04
$ws.Range("D39").Value = 11
$ws.Range("B40").Value = 3358
$ws.Range("C40").Value = 2.72
$ws.Range("E40").Value = 2346
$ws.Range("B41").Value = 378
$ws.Range("C41").Value = 1.1399999999999999
$ws.Range("D41").Value = 10.9
$ws.Range("E41").Value = 978
$ws.Range("F41").Value = 0.26
$ws.Range("B42").Value = 2400
$ws.Range("C42").Value = 3.07
$ws.Range("D42").Value = 9
$ws.Range("E42").Value = 960
$ws.Range("F42").Value = 0.28000000000000003
$ws.Range("B43").Value = 93
$ws.Range("C43").Value = 0.97
$ws.Range("D43").Value = 6.8
$ws.Range("E43").Value = 97
$ws.Range("F43").Value = 0.28000000000000003
$ws.Range("G43").Value = 4.8
$ws.Range("B44").Value = 487
$ws.Range("C44").Value = 2.6
$ws.Range("D44").Value = 8.1
$ws.Range("E44").Value = 312
$ws.Range("G44").Value = 4.9000000000000004
$ws.Range("C45").Value = 3.18
$ws.Range("D45").Value = 9.8000000000000007
$ws.Range("E45").Value = 4308
$ws.Range("F45").Value = 0.27
$ws.Range("G45").Value = 5.2
$ws.Range("H45").Value = 645
$ws.Range("I45").Value = 1.31
$ws.Range("J45").Value = 19.3
$ws.Range("E46").Value = 978
$ws.Range("F46").Value = 0.26
$ws.Range("G46").Value = 5.4
$ws.Range("C47").Value = 3.18
$ws.Range("D47").Value = 9.8000000000000007
$ws.Range("E47").Value = 155
$ws.Range("G47").Value = 4.9000000000000004
$ws.Range("H47").Value = 275
$ws.Range("I47").Value = 0.55000000000000004
$ws.Range("J47").Value = 14.8
$ws.Range("E48").Value = 743
$ws.Range("F48").Value = 0.26
$ws.Range("G48").Value = 4.8
$ws.Range("E49").Value = 2432
$ws.Range("F49").Value = 0.27
$ws.Range("H49").Value = 370
$ws.Range("I49").Value = 1.98
$ws.Range("J49").Value = 23.3
$ws.Range("B50").Value = 2423
$ws.Range("C50").Value = 0.62
$ws.Range("E50").Value = 4851
$ws.Range("F50").Value = 0.48
$ws.Range("G50").Value = 9.1
$ws.Range("B51").Value = 567
$ws.Range("C51").Value = 0.62
$ws.Range("D51").Value = 10.8
$ws.Range("E51").Value = 806
$ws.Range("F51").Value = 0.72
$ws.Range("G51").Value = 11.9
$ws.Range("B52").Value = 184
$ws.Range("D52").Value = 12.4
$ws.Range("E52").Value = 1332
$ws.Range("F52").Value = 0.31
$ws.Range("G52").Value = 5.5
$ws.Range("B55").Value = 36
$ws.Range("C55").Value = 0.42
$ws.Range("D55").Value = 6
$ws.Range("B56").Value = 515
$ws.Range("C56").Value = 0.76
$ws.Range("D56").Value = 21.9
$ws.Range("E56").Value = 557
$ws.Range("G56").Value = 22.1
$ws.Range("B57").Value = 1121
$ws.Range("C57").Value = 0.57999999999999996
$ws.Range("D57").Value = 13.1
$ws.Range("E57").Value = 76
$ws.Range("F57").Value = 1.07
$ws.Range("E58").Value = 2081
$ws.Range("F58").Value = 0.4
$ws.Range("G58").Value = 6.9
$ws.Range("E59").Value = 70
$ws.Range("E61").Value = 70
$ws.Range("H63").Value = 8
$ws.Range("I63").Value = 0.16
$ws.Range("J63").Value = 10.5
$ws.Range("H64").Value = 8
$ws.Range("I64").Value = 0.16
$ws.Range("J64").Value = 10.5
$ws.Range("B66").Value = 14908
$ws.Range("C66").Value = 2.2999999999999998
$ws.Range("D66").Value = 10.1
$ws.Range("E66").Value = 24028
$ws.Range("G66").Value = 6
$ws.Range("H66").Value = 2408
$ws.Range("I66").Value = 0.92
$ws.Range("J66").Value = 12.7
